# The document has a first-page footer/header that differ from the
# default (non-first-page) ones. Three inline pictures had their
# internal "name" bookkeeping values swapped:
#   - the Pearson logo living in the *default* footer keeps the name
#     "image1.png" in the source file, but should become "image2.png"
#   - the Pearson logo living in the *first-page* footer likewise
#     changes from "image1.png" to "image2.png"
#   - the BTEC logo living in the *first-page* header changes from
#     "image2.jpg" to "image1.jpg"
# (The alt-text/description stays untouched - only the shape's Name.)

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Default footer (wdHeaderFooterPrimary = 1): Pearson logo ---
$footerPrimary = $sec.Footers.Item(1)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    $pearsonPrimary = $footerPrimary.Range.InlineShapes.Item(1)
    $pearsonPrimary.Name = "image2.png"
}

# --- First-page footer (wdHeaderFooterFirstPage = 2): Pearson logo ---
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $pearsonFirst = $footerFirst.Range.InlineShapes.Item(1)
    $pearsonFirst.Name = "image2.png"
}

# --- First-page header (wdHeaderFooterFirstPage = 2): BTEC logo ---
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $btec = $headerFirst.Range.InlineShapes.Item(1)
    $btec.Name = "image1.jpg"
}

Write-Host "Renamed inline shapes complete"
